# Update the NATMI TPM-derived receptor/edge metrics in the sheet.
# The underlying "Receptor average expression value" (column M) for the
# target cluster "ECs" (rows 2 and 5) was recomputed with new TPM values,
# which cascades into the Receptor total expression value (N) and the
# derived specificity / edge-weight columns (O, P, Q, R, S, T) for every
# row, because those are normalized across the three target clusters
# (ECs, FAPs, MuSCs) present in rows 2-4 / 5-7.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "M2" = 20.56839166666667
    "N2" = 61.705175
    "O2" = 0.1304525281245593
    "P2" = 0.1304525281245593
    "Q2" = 0.4061640302416667
    "R2" = 3.655476272175
    "S2" = 0.0931662232504764
    "T2" = 0.09316622325047641

    "O3" = 0.6526310778549473
    "P3" = 0.6526310778549473
    "S3" = 0.4660942457288118
    "T3" = 0.4660942457288118

    "O4" = 0.2169163940204933
    "P4" = 0.2169163940204934
    "Q4" = 0.6753693323346667
    "R4" = 6.078323991012
    "S4" = 0.1549167462105852
    "T4" = 0.1549167462105853

    "M5" = 20.56839166666667
    "N5" = 61.705175
    "O5" = 0.1304525281245593
    "P5" = 0.1304525281245593
    "Q5" = 0.1625519993416667
    "R5" = 1.462967994075
    "S5" = 0.0372863048740829
    "T5" = 0.03728630487408291

    "O6" = 0.6526310778549473
    "P6" = 0.6526310778549473
    "S6" = 0.1865368321261356
    "T6" = 0.1865368321261356

    "O7" = 0.2169163940204933
    "P7" = 0.2169163940204934
    "S7" = 0.0619996478099081
    "T7" = 0.06199964780990812
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
